# Case_Data.xlsx edit
# - Clears the stray empty G195 cell (was an empty inline-string placeholder).
# - Appends 9 new case rows (198-206) for case 21TRD09437 covering Bunner and
#   Hemmeter, mirroring the existing Dismissed/Guilty row pattern. Some of the
#   new fine/cost cells keep the "$ " prefix and some don't (inconsistent
#   clean-up per the commit message), and a couple of cells are left as
#   single-space placeholders / an explicit blank cell, matching the source
#   data exactly.
# - The worksheet's used-range <dimension> is recalculated by Excel
#   automatically once the new cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 195: remove the stray empty G195 placeholder cell ---------------
$ws.Range("G195").ClearContents()

# --- Helper: write a value while forcing text storage for cells whose ----
# --- content would otherwise be auto-parsed as a number by Excel ---------
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# --- Row 198 ---------------------------------------------------------------
$ws.Range("A198").Value = "21TRD09437"
$ws.Range("B198").Value = "Bunner"
$ws.Range("C198").Value = "DUS"
Set-TextCell $ws "D198" "4510.11"
$ws.Range("E198").Value = "M1"
$ws.Range("F198").Value = "Dismissed"
$ws.Range("H198").Value = " "
$ws.Range("I198").Value = " "

# --- Row 199 ---------------------------------------------------------------
$ws.Range("A199").Value = "21TRD09437"
$ws.Range("B199").Value = "Bunner"
$ws.Range("C199").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Range("D199").Value = "4511.21B1A"
$ws.Range("E199").Value = "M4"
$ws.Range("F199").Value = "Guilty"
$ws.Range("G199").Value = "Guilty"
Set-TextCell $ws "H199" "$ 0"
Set-TextCell $ws "I199" "0"

# --- Row 200 ---------------------------------------------------------------
$ws.Range("A200").Value = "21TRD09437"
$ws.Range("B200").Value = "Bunner"
$ws.Range("C200").Value = "RECKLESS OPERATION 1ST IN 1 YR"
Set-TextCell $ws "D200" "4511.20"
$ws.Range("E200").Value = "MM"
$ws.Range("F200").Value = "Guilty"
$ws.Range("G200").Value = "Guilty"
Set-TextCell $ws "H200" "$ 0"
Set-TextCell $ws "I200" "0"

# --- Row 201 ---------------------------------------------------------------
$ws.Range("A201").Value = "21TRD09437"
$ws.Range("B201").Value = "Hemmeter"
$ws.Range("C201").Value = "DUS"
Set-TextCell $ws "D201" "4510.11"
$ws.Range("E201").Value = "M1"
$ws.Range("F201").Value = "Guilty"
$ws.Range("G201").Value = "Guilty"
Set-TextCell $ws "H201" "$ 25"
Set-TextCell $ws "I201" "$ 10"

# --- Row 202 ---------------------------------------------------------------
$ws.Range("A202").Value = "21TRD09437"
$ws.Range("B202").Value = "Hemmeter"
$ws.Range("C202").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Range("D202").Value = "4511.21B1A"
$ws.Range("E202").Value = "M4"
$ws.Range("F202").Value = "Dismissed"
$ws.Range("H202").Value = " "
$ws.Range("I202").Value = " "

# --- Row 203 ---------------------------------------------------------------
$ws.Range("A203").Value = "21TRD09437"
$ws.Range("B203").Value = "Hemmeter"
$ws.Range("C203").Value = "RECKLESS OPERATION 1ST IN 1 YR"
Set-TextCell $ws "D203" "4511.20"
$ws.Range("E203").Value = "MM"
$ws.Range("F203").Value = "Guilty"
$ws.Range("G203").Value = "Guilty"
Set-TextCell $ws "H203" "$ 0"
Set-TextCell $ws "I203" "$ 0"

# --- Row 204 ---------------------------------------------------------------
$ws.Range("A204").Value = "21TRD09437"
$ws.Range("B204").Value = "Bunner"
$ws.Range("C204").Value = "DUS"
Set-TextCell $ws "D204" "4510.11"
$ws.Range("E204").Value = "M1"
$ws.Range("F204").Value = "Guilty"
$ws.Range("G204").Value = "Guilty"
Set-TextCell $ws "H204" "$ 50"
Set-TextCell $ws "I204" "$ 25"
$ws.Range("J204").Value = "None"
$ws.Range("K204").Value = "None"

# --- Row 205 (G205 is an explicit blank cell, not absent) ------------------
$ws.Range("A205").Value = "21TRD09437"
$ws.Range("B205").Value = "Bunner"
$ws.Range("C205").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Range("D205").Value = "4511.21B1A"
$ws.Range("E205").Value = "M4"
$ws.Range("F205").Value = "Dismissed"
$ws.Range("G205").Font.Bold = $false
$ws.Range("H205").Value = " "
$ws.Range("I205").Value = " "
$ws.Range("J205").Value = " "
$ws.Range("K205").Value = " "

# --- Row 206 ---------------------------------------------------------------
$ws.Range("A206").Value = "21TRD09437"
$ws.Range("B206").Value = "Bunner"
$ws.Range("C206").Value = "RECKLESS OPERATION 1ST IN 1 YR"
Set-TextCell $ws "D206" "4511.20"
$ws.Range("E206").Value = "MM"
$ws.Range("F206").Value = "Guilty"
$ws.Range("G206").Value = "Guilty"
Set-TextCell $ws "H206" "$ 0"
Set-TextCell $ws "I206" "$ 0"
$ws.Range("J206").Value = "None"
$ws.Range("K206").Value = "None"
